$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update dSF (column F) values per repulled / recalculated data
$ws.Range("F9").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("F21").Value = 0
$ws.Range("F32").Value = 1
$ws.Range("F38").Value = -2
$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 1
$ws.Range("F44").Value = 1
$ws.Range("F56").Value = -1
$ws.Range("F61").Value = -2
$ws.Range("F66").Value = 0
$ws.Range("F70").Value = 0
